# 2.1.1.1e — add a new "2023" data column (Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column Q: year 2023 ------------------------------------------------

# Header cell (year number), same look as the other year headers in row 4.
$ws.Range("Q4").Value = 2023

# The data values for every indicator row.
$ws.Range("Q5").Value = 74.605426356589135
$ws.Range("Q6").Value = 118.8
$ws.Range("Q7").Value = 71.61643835616438
$ws.Range("Q8").Value = 95.703125
$ws.Range("Q9").Value = 113.91018619934282
$ws.Range("Q10").Value = 108.21501014198785
$ws.Range("Q11").Value = 165.26684164479443
$ws.Range("Q12").Value = 48.504446240905416
$ws.Range("Q13").Value = 97.361348644026393
$ws.Range("Q14").Value = 52.747252747252752

# Match formatting to the existing last column (P) for each of those rows,
# including the empty, bottom-bordered Q3 cell.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)

$ws.Range("P5:P13").Copy()
$ws.Range("Q5:Q13").PasteSpecial(-4122)

$ws.Range("P14").Copy()
$ws.Range("Q14").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row heights -------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 16.5
$ws.Rows.Item(5).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 16.5
$ws.Rows.Item(8).RowHeight = 16.5
$ws.Rows.Item(9).RowHeight = 16.5
$ws.Rows.Item(10).RowHeight = 16.5
$ws.Rows.Item(11).RowHeight = 16.5
$ws.Rows.Item(12).RowHeight = 16.5
$ws.Rows.Item(13).RowHeight = 16.5
$ws.Rows.Item(14).RowHeight = 16.5
